$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("B15").Value = "Beste klant,`nDank je wel voor je e-mail. Kun je alsjeblieft meer details geven over wat je precies geregeld wilt hebben? Op die manier kan ik je beter helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("C15").Value = "Kun jij dit even regelen?"
$ws.Range("D15").Value = "mailmind.test@zohomail.eu"
$ws.Range("E15").Value = "Overig"
$ws.Range("F15").Value = "2025-08-01 22:51:18"
$ws.Range("G15").Value = "Ja"
$ws.Range("H15").Value = "Nee"
$ws.Range("I15").Value = "Ja"
$ws.Range("J15").Value = "Nee"

$ws.Rows.Item(15).AutoFit() | Out-Null
